# Updated input for current test case
# Remove our design test values from the input spreadsheet and replace
# them with correct test values for Max's test case.

$wb = $excel.ActiveWorkbook

# --- Sources sheet: clear the "frozen/scrolled" top-left cell of the view ---
$wsSources = $wb.Worksheets.Item("Sources")
$wsSources.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# --- Transformers sheet: drop the "STEVE-O" test row and simplify the
#     Refinery row so it only yields a single product (gasoline) ---
$wsTrans = $wb.Worksheets.Item("Transformers")
$wsTrans.Activate()

# Refinery (row 2): SubEff0 becomes 1 (100%), and Prod1/SubEff1/Prod2/SubEff2
# (columns H:K) are no longer used.
$wsTrans.Range("G2").Value = 1
$wsTrans.Range("H2:K2").ClearContents()

# STEVE-O (row 4) is removed entirely - delete the whole row.
$wsTrans.Rows("4:4").Delete()

# Update the selected cell shown when the sheet is opened.
$wsTrans.Range("G13").Select() | Out-Null

$wb.Save() | Out-Null
